$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Total" row (row 6) mirroring the style of the preceding data rows
$ws.Range("A6").Value = "Total"
$ws.Range("B6").Value = 6003
$ws.Range("C6").Value = 1489
$ws.Range("D6").Value = 7492
$ws.Range("E6").Value = ""

# Match styling used by the other data rows (A2:E5): centered alignment, no border/bold
$ws.Range("A6:E6").HorizontalAlignment = -4108
